$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, in the order cells appear in the sheet.
# Column D holds price strings that look numeric (e.g. "1.003"); Excel would
# auto-convert those to numbers on plain assignment, so for "D" cells we force
# the Text number format before assigning, then restore the default "Normal"
# style afterwards so the saved cell keeps no special formatting (matching the
# original workbook) while the stored value remains the exact text string.
$updates = [ordered]@{
    'D2' = '30.281.92'
    'E2' = '  +1.86%  '
    'D3' = '2.088.49'
    'E3' = '  -0.48%  '
    'D4' = '1.003'
    'E4' = '  -0.49%  '
    'D5' = '342.68'
    'E5' = '  -0.69%  '
    'E6' = '  -0.35%  '
    'D7' = '0.5228'
    'E7' = '  +1.70%  '
    'D8' = '0.4403'
    'E8' = '  -0.16%  '
    'D9' = '54.38'
    'E9' = '  +3.09%  '
    'D10' = '0.09350'
    'E10' = '  +1.17%  '
    'D11' = '1.167'
    'E11' = '  -0.44%  '
    'D12' = '24.76'
    'E12' = '  -0.58%  '
    'D13' = '8.568'
    'E13' = '  +3.41%  '
    'D14' = '6.889'
    'E14' = '  +1.92%  '
    'D15' = '2.086.18'
    'E15' = '  -0.70%  '
    'D16' = '101.46'
    'E16' = '  +1.89%  '
    'E17' = '  +0.42%  '
    'D18' = '1.004'
    'E18' = '  -0.40%  '
    'D19' = '21.08'
    'E19' = '  +1.08%  '
    'D20' = '0.06668'
    'E20' = '  +0.59%  '
    'D21' = '6.319'
    'E21' = '  +1.98%  '
    'E22' = '  -0.36%  '
    'D23' = '30.285.34'
    'E23' = '  +1.66%  '
    'D24' = '12.50'
    'E24' = '  -0.95%  '
    'D25' = '2.308'
    'E25' = '  -0.50%  '
    'D26' = '21.76'
    'E26' = '  -0.71%  '
    'D27' = '162.13'
    'E27' = '  +0.12%  '
    'D28' = '2.504'
    'E28' = '  -1.02%  '
    'D29' = '132.88'
    'E29' = '  -0.05%  '
    'D30' = '1.130'
    'E30' = '  -0.35%  '
    'B31' = 'ARBITRUM'
    'C31' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D31' = '1.659'
    'E31' = '  +0.09%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.1044'
    'E32' = '  -0.62%  '
    'D33' = '6.211'
    'E33' = '  +0.46%  '
    'D34' = '6.648'
    'E34' = '  +9.92%  '
    'D35' = '3.875'
    'E35' = '  -1.59%  '
    'D36' = '10.21'
    'E36' = '  -2.72%  '
    'D37' = '0.02622'
    'E37' = '  +2.01%  '
    'D38' = '0.06817'
    'E38' = '  +1.18%  '
    'D39' = '0.6962'
    'E39' = '  +1.39%  '
    'D40' = '1.339'
    'E40' = '  +2.82%  '
    'D41' = '12.50'
    'E41' = '  +0.33%  '
    'E42' = '  -1.24%  '
    'D43' = '0.6799'
    'E43' = '  +2.27%  '
    'D44' = '14.23'
    'E44' = '  -0.04%  '
    'D45' = '2.322'
    'E45' = '  +0.21%  '
    'E46' = '  -0.27%  '
    'E47' = '  +17.80%  '
    'D48' = '3.630'
    'E48' = '  +0.25%  '
    'E49' = '  +0.41%  '
    'D50' = '1.209'
    'E50' = '  +7.94%  '
    'E51' = '  -0.68%  '
}

foreach ($cellRef in $updates.Keys) {
    $value = $updates[$cellRef]
    $cell = $ws.Range($cellRef)
    if ($cellRef -match "^D\d+$") {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
